# Set Exis Unit (column F) to 0 and MaxInvest (column I) to 200 ("MaxlineLoad 100%")
# for data rows 8-18 on the active sheet, then update the active selection
# to reflect the F9:F18 range (matching the edited workbook's selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 8; $row -le 18; $row++) {
    $ws.Cells.Item($row, 6).Value = 0     # Column F - ExisUnits
    $ws.Cells.Item($row, 9).Value = 200   # Column I - MaxInvest
}

$ws.Range("F9:F18").Select()
